# Apply the cryptos list update (values sourced from the coinranking feed).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.988.75"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.33%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.543.37"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -2.52%  "

# Row 4
$ws.Range("E4").Value = "  -0.12%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "199.13"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.92%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "558.19"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.98%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.649"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +5.15%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.535.38"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.13%  "

# Row 9
$ws.Range("E9").Value = "  +0.02%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.666"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.59%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "61.57"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +6.79%  "

# Row 12
$ws.Range("E12").Value = "  -5.50%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000271"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -6.60%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.02"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.22%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.100.83"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.90%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.539.45"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.76%  "

# Row 17
$ws.Range("E17").Value = "  -1.61%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.829.39"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.49%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "18.50"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.34%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.96"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -4.08%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.04"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -4.46%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "399.99"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.28%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.02"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -4.03%  "

# Row 24
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "86.03"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.39%  "

# Row 25
$ws.Range("B25").Value = "RenderToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.87"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -8.23%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.49"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.71%  "

# Row 27
$ws.Range("B27").Value = "ImmutableX"
$ws.Range("C27").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.85"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -3.44%  "

# Row 28
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.87"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.71%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.96"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.04%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "724.86"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +3.30%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "31.41"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.30%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.13"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -12.96%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.83"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.29%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "64.34"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.27%  "

# Row 35
$ws.Range("E35").Value = "  -3.40%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "38.85"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -8.83%  "

# Row 38
$ws.Range("E38").Value = "  -6.79%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.133"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -5.30%  "

# Row 40
$ws.Range("E40").Value = "  -3.13%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.091.86"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -4.24%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.998"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.21%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0₃0691"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -11.54%  "

# Row 44
$ws.Range("E44").Value = "  -10.04%  "

# Row 45
$ws.Range("E45").Value = "  +1.87%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0413"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.68%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.135"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.76%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.60"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -13.11%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "139.79"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.50%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.02"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.12%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.31"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -6.77%  "
